$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 933.875
$ws.Range("I28").Value = 133.35715
$ws.Range("J28").Value = 2054.6
$ws.Range("K28").Value = 133.35715
$ws.Range("L28").Value = 2054.6
$ws.Range("M28").Value = 351.64285
$ws.Range("N28").Value = -3024.6
# Row 116
$ws.Range("H116").Value = 2239.6
$ws.Range("I116").Value = 2249.5
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 2249.5
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = 1192.5
$ws.Range("N116").Value = -9084
# Row 134
$ws.Range("H134").Value = 29756
$ws.Range("J134").Value = 29756
$ws.Range("L134").Value = 29756
$ws.Range("N134").Value = -39896

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 834.9
$ws.Range("I45").Value = 796.8
$ws.Range("K45").Value = 796.8
$ws.Range("M45").Value = -419.8
# Row 61
$ws.Range("H61").Value = 1772.9565
$ws.Range("I61").Value = 1106.8
$ws.Range("K61").Value = 1106.8
$ws.Range("M61").Value = -894.8
# Row 74
$ws.Range("H74").Value = 799.2286
$ws.Range("I74").Value = 761.8929000000001
$ws.Range("J74").Value = 948.5714
$ws.Range("K74").Value = 761.8929000000001
$ws.Range("L74").Value = 948.5714
$ws.Range("M74").Value = 112.1070999999999
$ws.Range("N74").Value = -2696.5714
# Row 77
$ws.Range("H77").Value = 799.2286
$ws.Range("I77").Value = 761.8929000000001
$ws.Range("J77").Value = 948.5714
$ws.Range("K77").Value = 3809.4645
$ws.Range("L77").Value = 4742.857
$ws.Range("M77").Value = 558.5355
$ws.Range("N77").Value = -13478.857
# Row 132
$ws.Range("H132").Value = 4756.9033
$ws.Range("I132").Value = 5857.5
$ws.Range("J132").Value = 2066.5557
$ws.Range("K132").Value = 17572.5
$ws.Range("L132").Value = 6199.6671
$ws.Range("M132").Value = -15042.5
$ws.Range("N132").Value = -11259.6671
# Row 134
$ws.Range("H134").Value = 51000
$ws.Range("J134").Value = 51000
$ws.Range("L134").Value = 51000
$ws.Range("N134").Value = -61140
# Row 136
$ws.Range("H136").Value = 1772.9565
$ws.Range("I136").Value = 1106.8
$ws.Range("K136").Value = 3320.4
$ws.Range("M136").Value = -770.3999999999996

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3802.6667
$ws.Range("I105").Value = 3498.75
$ws.Range("K105").Value = 3498.75
$ws.Range("M105").Value = -1751.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9092727
$ws.Range("I31").Value = 2250.7144
$ws.Range("J31").Value = 25001060
$ws.Range("K31").Value = 2250.7144
$ws.Range("L31").Value = 25001060
$ws.Range("M31").Value = -1955.7144
$ws.Range("N31").Value = -25001650
# Row 34
$ws.Range("H34").Value = 9092727
$ws.Range("I34").Value = 2250.7144
$ws.Range("J34").Value = 25001060
$ws.Range("K34").Value = 2250.7144
$ws.Range("L34").Value = 25001060
$ws.Range("M34").Value = -2048.7144
$ws.Range("N34").Value = -25001464
# Row 99
$ws.Range("H99").Value = 1807.5667
$ws.Range("I99").Value = 1501.35
$ws.Range("J99").Value = 2420
$ws.Range("K99").Value = 1501.35
$ws.Range("L99").Value = 2420
$ws.Range("M99").Value = -3.349999999999909
$ws.Range("N99").Value = -5416
# Row 126
$ws.Range("H126").Value = 1807.5667
$ws.Range("I126").Value = 1501.35
$ws.Range("J126").Value = 2420
$ws.Range("K126").Value = 4504.049999999999
$ws.Range("L126").Value = 7260
$ws.Range("M126").Value = -2034.049999999999
$ws.Range("N126").Value = -12200

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 614.4
$ws.Range("I5").Value = 156.66667
$ws.Range("J5").Value = 1301
$ws.Range("K5").Value = 470.00001
$ws.Range("L5").Value = 3903
$ws.Range("M5").Value = -358.00001
$ws.Range("N5").Value = -4127
# Row 131
$ws.Range("H131").Value = 6179709.5
$ws.Range("J131").Value = 7937343.5
$ws.Range("L131").Value = 23812030.5
$ws.Range("N131").Value = -23822110.5
# Row 135
$ws.Range("H135").Value = 614.4
$ws.Range("I135").Value = 156.66667
$ws.Range("J135").Value = 1301
$ws.Range("K135").Value = 1410.00003
$ws.Range("L135").Value = 11709
$ws.Range("M135").Value = 1124.99997
$ws.Range("N135").Value = -16779

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 539.069
$ws.Range("I107").Value = 488.09525
$ws.Range("J107").Value = 672.875
$ws.Range("K107").Value = 488.09525
$ws.Range("L107").Value = 672.875
$ws.Range("M107").Value = 1431.90475
$ws.Range("N107").Value = -4512.875
# Row 122
$ws.Range("H122").Value = 2157.1428
$ws.Range("I122").Value = 1847.0588
$ws.Range("J122").Value = 2636.3635
$ws.Range("K122").Value = 5541.1764
$ws.Range("L122").Value = 7909.0905
$ws.Range("M122").Value = -3091.1764
$ws.Range("N122").Value = -12809.0905
# Row 126
$ws.Range("H126").Value = 1264.7059
$ws.Range("J126").Value = 1750
$ws.Range("L126").Value = 5250
$ws.Range("N126").Value = -10190

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 9782.909
$ws.Range("I16").Value = 701.5
$ws.Range("J16").Value = 34000
$ws.Range("K16").Value = 701.5
$ws.Range("L16").Value = 34000
$ws.Range("M16").Value = -531.5
$ws.Range("N16").Value = -34340
# Row 46
$ws.Range("H46").Value = 1871.0869
$ws.Range("I46").Value = 1583.5
$ws.Range("J46").Value = 1972.5883
$ws.Range("K46").Value = 1583.5
$ws.Range("L46").Value = 1972.5883
$ws.Range("M46").Value = -1395.5
$ws.Range("N46").Value = -2348.5883
# Row 55
$ws.Range("H55").Value = 277.75
$ws.Range("I55").Value = 285.7143
$ws.Range("J55").Value = 222
$ws.Range("K55").Value = 285.7143
$ws.Range("L55").Value = 222
$ws.Range("M55").Value = -112.7143
$ws.Range("N55").Value = -568

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1699.4
$ws.Range("I96").Value = 1699.4
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1699.4
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -326.4000000000001
$ws.Range("N96").ClearContents()
# Row 113
$ws.Range("H113").Value = 527.93335
$ws.Range("I113").Value = 527.93335
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1583.80005
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 586.1999499999999
$ws.Range("N113").ClearContents()
# Row 122
$ws.Range("H122").Value = 2011.1111
$ws.Range("I122").Value = 2233.3333
$ws.Range("J122").Value = 1566.6666
$ws.Range("K122").Value = 6699.999899999999
$ws.Range("L122").Value = 4699.9998
$ws.Range("M122").Value = -4249.999899999999
$ws.Range("N122").Value = -9599.9998
# Row 132
$ws.Range("H132").Value = 1225.0938
$ws.Range("I132").Value = 1097.1072
$ws.Range("J132").Value = 2121
$ws.Range("K132").Value = 3291.3216
$ws.Range("L132").Value = 6363
$ws.Range("M132").Value = -761.3215999999998
$ws.Range("N132").Value = -11423
# Row 136
$ws.Range("H136").Value = 7193.55
$ws.Range("I136").Value = 8257.117
$ws.Range("J136").Value = 1166.6666
$ws.Range("K136").Value = 24771.351
$ws.Range("L136").Value = 3499.9998
$ws.Range("M136").Value = -22221.351
$ws.Range("N136").Value = -8599.9998

